# Fix a row-ordering bug introduced by an earlier merge: the data rows
# (2-23, columns A:F) ended up permuted incorrectly. Restore the correct
# row order/values by writing the corrected table back out in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(101,  9, 30, 15, 60, 15),
    @(401,  9, 48, 67, 75, 45),
    @(902,  1,  0,  0,  0,  0),
    @(701,  3, 90, 45, 97, 15),
    @(801,  3, 67, 65, 52, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(1001,18, 30, 75, 60, 72),
    @(301,  6, 45, 30, 60, 45),
    @(501,  9, 52, 30, 75, 45),
    @(601,  9, 60, 67, 60, 42),
    @(201,  9, 30, 15, 45, 30),
    @(1201, 2, 10, 10, 10, 10),
    @(901, 16, 15, 45, 60, 60),
    @(1203, 3, 15, 15, 15, 15),
    @(1101, 0, 15, 30, 30,  0),
    @(802,  0,  4,  5,  4,  0),
    @(502,  0,  4,  0,  0,  0),
    @(2,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(1,    0,  2,  2,  2,  2),
    @(402,  0,  0,  4,  0,  0),
    @(602,  0,  0,  4,  0,  9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
